$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update window position (workbookView xWindow/yWindow)
$excel.Left = 2580
$excel.Top = 660

# Update formulas/values
$ws.Range("B5").Formula = "=5+1"
$ws.Range("B6").Formula = "=4"

# Update selected cell
$ws.Range("B7").Select()

$wb.Save()
